$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.158.12"
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = "1.903.32"
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = "'306.05"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = "'0.5234"
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").Value = "'0.3761"
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("D9").Value = "'0.07253"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").Value = "'21.12"
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").Value = "'0.9020"
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").Value = "'0.08471"
$ws.Range("E12").Value = '  +11.11%  '
$ws.Range("D13").Value = "1.889.59"
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").Value = "'94.99"
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = "'5.289"
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = "'0.000008649"
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("D18").Value = "'14.54"
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").Value = "27.197.19"
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").Value = "'5.068"
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = "2.141.09"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = "'6.425"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = "'147.47"
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").Value = "'2.278"
$ws.Range("E26").Value = '  +3.75%  '
$ws.Range("D27").Value = "'1.753"
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").Value = "'18.19"
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = "'114.94"
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").Value = "'4.811"
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("D31").Value = "'4.887"
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("D32").Value = "'0.09261"
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("D33").Value = "'0.8101"
$ws.Range("E33").Value = '  +5.05%  '
$ws.Range("D34").Value = "'0.05058"
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("D35").Value = "'1.234"
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").Value = "'3.437"
$ws.Range("E36").Value = '  +4.64%  '
$ws.Range("D37").Value = "'2.945"
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("D38").Value = "'2.620"
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").Value = "'0.5694"
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("D40").Value = "'0.01994"
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").Value = "'1.076"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'9.006"
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").Value = "'6.626"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = "'116.68"
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = "'0.1514"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = "'0.4855"
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = "'10.15"
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").Value = "'1.615"
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").Value = "'37.50"
$ws.Range("E50").Value = '  +0.18%  '
$ws.Range("D51").Value = "'63.98"
$ws.Range("E51").Value = '  -0.52%  '
